$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New survey respondents: rows 188-223 ---
$ws.Cells.Item(188, 1).Value = "Esther"; $ws.Cells.Item(188, 2).Value = 36; $ws.Cells.Item(188, 3).Value = "F"; $ws.Cells.Item(188, 4).Value = 14; $ws.Cells.Item(188, 6).Value = 66000; $ws.Cells.Item(188, 7).Value = 0; $ws.Cells.Item(188, 8).Value = -1; $ws.Cells.Item(188, 9).Value = 1
$ws.Cells.Item(189, 1).Value = "Nala"; $ws.Cells.Item(189, 2).Value = 20; $ws.Cells.Item(189, 3).Value = "F"; $ws.Cells.Item(189, 4).Value = 30; $ws.Cells.Item(189, 6).Value = 0; $ws.Cells.Item(189, 7).Value = 0; $ws.Cells.Item(189, 8).Value = -1; $ws.Cells.Item(189, 9).Value = 0
$ws.Cells.Item(190, 1).Value = "May"; $ws.Cells.Item(190, 2).Value = 26; $ws.Cells.Item(190, 3).Value = "F"; $ws.Cells.Item(190, 4).Value = 14; $ws.Cells.Item(190, 6).Value = 38000; $ws.Cells.Item(190, 7).Value = 0; $ws.Cells.Item(190, 8).Value = -1; $ws.Cells.Item(190, 9).Value = 0
$ws.Cells.Item(191, 1).Value = "June"; $ws.Cells.Item(191, 2).Value = 28; $ws.Cells.Item(191, 3).Value = "F"; $ws.Cells.Item(191, 4).Value = 11; $ws.Cells.Item(191, 6).Value = -1; $ws.Cells.Item(191, 7).Value = 0; $ws.Cells.Item(191, 8).Value = -1; $ws.Cells.Item(191, 9).Value = 1
$ws.Cells.Item(192, 1).Value = "Leo"; $ws.Cells.Item(192, 2).Value = 37; $ws.Cells.Item(192, 3).Value = "M"; $ws.Cells.Item(192, 4).Value = 5; $ws.Cells.Item(192, 6).Value = 67000; $ws.Cells.Item(192, 7).Value = 0; $ws.Cells.Item(192, 8).Value = 800; $ws.Cells.Item(192, 9).Value = 1
$ws.Cells.Item(193, 1).Value = "Phillip"; $ws.Cells.Item(193, 2).Value = 27; $ws.Cells.Item(193, 3).Value = "M"; $ws.Cells.Item(193, 4).Value = 5; $ws.Cells.Item(193, 6).Value = 46500; $ws.Cells.Item(193, 7).Value = 0; $ws.Cells.Item(193, 8).Value = -1; $ws.Cells.Item(193, 9).Value = 1
$ws.Cells.Item(194, 1).Value = "Yuri"; $ws.Cells.Item(194, 2).Value = 28; $ws.Cells.Item(194, 3).Value = "M"; $ws.Cells.Item(194, 4).Value = 11; $ws.Cells.Item(194, 6).Value = 57000; $ws.Cells.Item(194, 7).Value = 0; $ws.Cells.Item(194, 8).Value = -1; $ws.Cells.Item(194, 9).Value = 1; $ws.Cells.Item(194, 10).Value = "Yes"
$ws.Cells.Item(195, 1).Value = "Ernest"; $ws.Cells.Item(195, 2).Value = 24; $ws.Cells.Item(195, 3).Value = "M"; $ws.Cells.Item(195, 4).Value = 3; $ws.Cells.Item(195, 6).Value = 0; $ws.Cells.Item(195, 7).Value = -100; $ws.Cells.Item(195, 8).Value = -1; $ws.Cells.Item(195, 9).Value = 0
$ws.Cells.Item(196, 1).Value = "Ivy"; $ws.Cells.Item(196, 2).Value = 33; $ws.Cells.Item(196, 3).Value = "F"; $ws.Cells.Item(196, 4).Value = 30; $ws.Cells.Item(196, 6).Value = 71500; $ws.Cells.Item(196, 7).Value = -150; $ws.Cells.Item(196, 8).Value = 1200; $ws.Cells.Item(196, 9).Value = 1
$ws.Cells.Item(197, 1).Value = "Jonas"; $ws.Cells.Item(197, 2).Value = 28; $ws.Cells.Item(197, 3).Value = "M"; $ws.Cells.Item(197, 4).Value = 2; $ws.Cells.Item(197, 6).Value = 52000; $ws.Cells.Item(197, 7).Value = 0; $ws.Cells.Item(197, 8).Value = -1; $ws.Cells.Item(197, 9).Value = 1
$ws.Cells.Item(198, 1).Value = "Jason"; $ws.Cells.Item(198, 2).Value = 27; $ws.Cells.Item(198, 3).Value = "M"; $ws.Cells.Item(198, 4).Value = 3; $ws.Cells.Item(198, 6).Value = 48000; $ws.Cells.Item(198, 7).Value = 0; $ws.Cells.Item(198, 8).Value = -1; $ws.Cells.Item(198, 9).Value = 1
$ws.Cells.Item(199, 1).Value = "Phil"; $ws.Cells.Item(199, 2).Value = 24; $ws.Cells.Item(199, 3).Value = "M"; $ws.Cells.Item(199, 4).Value = 4; $ws.Cells.Item(199, 6).Value = 0; $ws.Cells.Item(199, 7).Value = 0; $ws.Cells.Item(199, 8).Value = -1; $ws.Cells.Item(199, 9).Value = 0
$ws.Cells.Item(200, 1).Value = "Kendrick"; $ws.Cells.Item(200, 2).Value = 25; $ws.Cells.Item(200, 3).Value = "F"; $ws.Cells.Item(200, 4).Value = 5; $ws.Cells.Item(200, 6).Value = 37000; $ws.Cells.Item(200, 7).Value = 0; $ws.Cells.Item(200, 8).Value = -1; $ws.Cells.Item(200, 9).Value = 0
$ws.Cells.Item(201, 1).Value = "Farah"; $ws.Cells.Item(201, 2).Value = 27; $ws.Cells.Item(201, 3).Value = "F"; $ws.Cells.Item(201, 4).Value = 12; $ws.Cells.Item(201, 6).Value = 51000; $ws.Cells.Item(201, 7).Value = 0; $ws.Cells.Item(201, 8).Value = -1; $ws.Cells.Item(201, 9).Value = 1
$ws.Cells.Item(202, 1).Value = "Lucy"; $ws.Cells.Item(202, 2).Value = 32; $ws.Cells.Item(202, 3).Value = "F"; $ws.Cells.Item(202, 4).Value = 14; $ws.Cells.Item(202, 6).Value = 62000; $ws.Cells.Item(202, 7).Value = 0; $ws.Cells.Item(202, 8).Value = -1; $ws.Cells.Item(202, 9).Value = 1
$ws.Cells.Item(203, 1).Value = "Lori"; $ws.Cells.Item(203, 2).Value = 38; $ws.Cells.Item(203, 3).Value = "F"; $ws.Cells.Item(203, 4).Value = 1; $ws.Cells.Item(203, 6).Value = 87000; $ws.Cells.Item(203, 7).Value = 0; $ws.Cells.Item(203, 8).Value = 2000; $ws.Cells.Item(203, 9).Value = 1
$ws.Cells.Item(204, 1).Value = "Scott"; $ws.Cells.Item(204, 2).Value = 38; $ws.Cells.Item(204, 3).Value = "M"; $ws.Cells.Item(204, 4).Value = 2; $ws.Cells.Item(204, 6).Value = 93000; $ws.Cells.Item(204, 7).Value = 0; $ws.Cells.Item(204, 8).Value = 2100; $ws.Cells.Item(204, 9).Value = 1
$ws.Cells.Item(205, 1).Value = "Waldo"; $ws.Cells.Item(205, 2).Value = 37; $ws.Cells.Item(205, 3).Value = "M"; $ws.Cells.Item(205, 4).Value = 4; $ws.Cells.Item(205, 5).Value = 1; $ws.Cells.Item(205, 6).Value = 87000; $ws.Cells.Item(205, 7).Value = -650; $ws.Cells.Item(205, 8).Value = 1800; $ws.Cells.Item(205, 9).Value = 1; $ws.Cells.Item(205, 11).Value = 1
$ws.Cells.Item(206, 1).Value = "Tomas"; $ws.Cells.Item(206, 2).Value = 28; $ws.Cells.Item(206, 3).Value = "M"; $ws.Cells.Item(206, 4).Value = 20; $ws.Cells.Item(206, 6).Value = 0; $ws.Cells.Item(206, 7).Value = 300; $ws.Cells.Item(206, 8).Value = -1; $ws.Cells.Item(206, 9).Value = -1
$ws.Cells.Item(207, 1).Value = "Miles"; $ws.Cells.Item(207, 2).Value = 25; $ws.Cells.Item(207, 3).Value = "M"; $ws.Cells.Item(207, 4).Value = 20; $ws.Cells.Item(207, 6).Value = 56000; $ws.Cells.Item(207, 7).Value = 0; $ws.Cells.Item(207, 8).Value = -1; $ws.Cells.Item(207, 9).Value = 1
$ws.Cells.Item(208, 1).Value = "Hope"; $ws.Cells.Item(208, 2).Value = 26; $ws.Cells.Item(208, 3).Value = "F"; $ws.Cells.Item(208, 4).Value = 20; $ws.Cells.Item(208, 6).Value = 28500; $ws.Cells.Item(208, 7).Value = 0; $ws.Cells.Item(208, 8).Value = -1; $ws.Cells.Item(208, 9).Value = -1
$ws.Cells.Item(209, 1).Value = "Ethan"; $ws.Cells.Item(209, 2).Value = 27; $ws.Cells.Item(209, 3).Value = "M"; $ws.Cells.Item(209, 4).Value = 13; $ws.Cells.Item(209, 6).Value = 41000; $ws.Cells.Item(209, 7).Value = 0; $ws.Cells.Item(209, 8).Value = -1; $ws.Cells.Item(209, 9).Value = 1
$ws.Cells.Item(210, 1).Value = "Abraham"; $ws.Cells.Item(210, 2).Value = 29; $ws.Cells.Item(210, 3).Value = "M"; $ws.Cells.Item(210, 4).Value = 21; $ws.Cells.Item(210, 6).Value = 74500; $ws.Cells.Item(210, 7).Value = -1000; $ws.Cells.Item(210, 8).Value = 1400; $ws.Cells.Item(210, 9).Value = 1
$ws.Cells.Item(211, 1).Value = "Joshua"; $ws.Cells.Item(211, 2).Value = 33; $ws.Cells.Item(211, 3).Value = "M"; $ws.Cells.Item(211, 4).Value = 13; $ws.Cells.Item(211, 6).Value = 67000; $ws.Cells.Item(211, 7).Value = -500; $ws.Cells.Item(211, 8).Value = 900; $ws.Cells.Item(211, 9).Value = 1
$ws.Cells.Item(212, 1).Value = "Abel"; $ws.Cells.Item(212, 2).Value = 34; $ws.Cells.Item(212, 3).Value = "M"; $ws.Cells.Item(212, 4).Value = 20; $ws.Cells.Item(212, 6).Value = 63500; $ws.Cells.Item(212, 7).Value = -400; $ws.Cells.Item(212, 8).Value = 850; $ws.Cells.Item(212, 9).Value = -2
$ws.Cells.Item(213, 1).Value = "Muhammad"; $ws.Cells.Item(213, 2).Value = 36; $ws.Cells.Item(213, 3).Value = "M"; $ws.Cells.Item(213, 4).Value = 1; $ws.Cells.Item(213, 6).Value = 79500; $ws.Cells.Item(213, 7).Value = -200; $ws.Cells.Item(213, 8).Value = 1200; $ws.Cells.Item(213, 9).Value = 0
$ws.Cells.Item(214, 1).Value = "Maged"; $ws.Cells.Item(214, 2).Value = 37; $ws.Cells.Item(214, 3).Value = "M"; $ws.Cells.Item(214, 4).Value = 1; $ws.Cells.Item(214, 6).Value = 77000; $ws.Cells.Item(214, 7).Value = -350; $ws.Cells.Item(214, 8).Value = 1300; $ws.Cells.Item(214, 9).Value = 0; $ws.Cells.Item(214, 10).Value = "Yes"
$ws.Cells.Item(215, 1).Value = "Cage"; $ws.Cells.Item(215, 2).Value = 31; $ws.Cells.Item(215, 3).Value = "M"; $ws.Cells.Item(215, 4).Value = 4; $ws.Cells.Item(215, 6).Value = 67000; $ws.Cells.Item(215, 7).Value = 0; $ws.Cells.Item(215, 8).Value = 800; $ws.Cells.Item(215, 9).Value = 1
$ws.Cells.Item(216, 1).Value = "Chase"; $ws.Cells.Item(216, 2).Value = 33; $ws.Cells.Item(216, 3).Value = "M"; $ws.Cells.Item(216, 4).Value = 3; $ws.Cells.Item(216, 6).Value = 59000; $ws.Cells.Item(216, 7).Value = -200; $ws.Cells.Item(216, 8).Value = -1; $ws.Cells.Item(216, 9).Value = 1; $ws.Cells.Item(216, 11).Value = 1
$ws.Cells.Item(217, 1).Value = "Victor"; $ws.Cells.Item(217, 2).Value = 28; $ws.Cells.Item(217, 3).Value = "M"; $ws.Cells.Item(217, 4).Value = 5; $ws.Cells.Item(217, 6).Value = 65000; $ws.Cells.Item(217, 7).Value = -100; $ws.Cells.Item(217, 8).Value = -1; $ws.Cells.Item(217, 9).Value = 0
$ws.Cells.Item(218, 1).Value = "Victoria"; $ws.Cells.Item(218, 2).Value = 31; $ws.Cells.Item(218, 3).Value = "F"; $ws.Cells.Item(218, 4).Value = 2; $ws.Cells.Item(218, 6).Value = 64500; $ws.Cells.Item(218, 7).Value = 100; $ws.Cells.Item(218, 8).Value = -1; $ws.Cells.Item(218, 9).Value = 0
$ws.Cells.Item(219, 1).Value = "Terrance"; $ws.Cells.Item(219, 2).Value = 28; $ws.Cells.Item(219, 3).Value = "M"; $ws.Cells.Item(219, 4).Value = 1; $ws.Cells.Item(219, 6).Value = 57500; $ws.Cells.Item(219, 7).Value = 0; $ws.Cells.Item(219, 8).Value = -1; $ws.Cells.Item(219, 9).Value = 0
$ws.Cells.Item(220, 1).Value = "Terri"; $ws.Cells.Item(220, 2).Value = 29; $ws.Cells.Item(220, 3).Value = "F"; $ws.Cells.Item(220, 4).Value = 20; $ws.Cells.Item(220, 6).Value = 9999999; $ws.Cells.Item(220, 7).Value = 0; $ws.Cells.Item(220, 8).Value = -1; $ws.Cells.Item(220, 9).Value = 1
$ws.Cells.Item(221, 1).Value = "Beth"; $ws.Cells.Item(221, 2).Value = 28; $ws.Cells.Item(221, 3).Value = "F"; $ws.Cells.Item(221, 4).Value = 30; $ws.Cells.Item(221, 6).Value = 73000; $ws.Cells.Item(221, 7).Value = 0; $ws.Cells.Item(221, 8).Value = -1; $ws.Cells.Item(221, 9).Value = 1
$ws.Cells.Item(222, 1).Value = "Ellen"; $ws.Cells.Item(222, 2).Value = 27; $ws.Cells.Item(222, 3).Value = "F"; $ws.Cells.Item(222, 4).Value = 30; $ws.Cells.Item(222, 6).Value = 71000; $ws.Cells.Item(222, 7).Value = 0; $ws.Cells.Item(222, 8).Value = -1; $ws.Cells.Item(222, 9).Value = 1
$ws.Cells.Item(223, 1).Value = "Elvis"; $ws.Cells.Item(223, 2).Value = 31; $ws.Cells.Item(223, 3).Value = "M"; $ws.Cells.Item(223, 4).Value = 20; $ws.Cells.Item(223, 6).Value = 68500; $ws.Cells.Item(223, 7).Value = 0; $ws.Cells.Item(223, 8).Value = -1; $ws.Cells.Item(223, 9).Value = 1

# --- Felonies (column K) backfilled for a few earlier respondents ---
$ws.Cells.Item(36, 11).Value = 2
$ws.Cells.Item(121, 11).Value = 1
$ws.Cells.Item(133, 11).Value = 2

# --- Move selection to the next empty row, as left by the author ---
$ws.Range("A224").Select()
